# LogBook.xlsx -- "Excel Integration To Measurement"
#
# Rows 2-6 of the log sheet get new/updated measurement entries and six
# brand-new rows (7-12) are appended, growing the used range from
# A1:N6 to A1:N12. Numeric-looking log values (frequencies, amplitudes,
# etc.) are written with a leading apostrophe so Excel stores them as
# literal text ("15.0") instead of silently re-typing them as numbers
# (15), matching how the log file itself represents them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = "None"
$ws.Range("B2").Value = "FieldSweep_27-10-22_16-42-35.dat"
$ws.Range("F2").Value = "150/250/0.4"

# --- Row 3 ---
$ws.Range("A3").Value = "None"
$ws.Range("B3").Value = "FieldSweep_27-10-22_16-36-43.dat"
$ws.Range("F3").Value = "150/250/0.4"
$ws.Range("G3").Value = "'468.0"
$ws.Range("K3").Value = "Repositioned Sample"

# --- Row 4 ---
$ws.Range("A4").Value = "None"
$ws.Range("B4").Value = "FieldSweep_27-10-22_15-34-48.dat"
$ws.Range("C4").Value = "None"
$ws.Range("D4").Value = "'15.0"
$ws.Range("E4").Value = "'13.0"
$ws.Range("F4").Value = "150/250/0.4"
$ws.Range("G4").Value = "'468.0"
$ws.Range("H4").Value = "'5.0"
$ws.Range("I4").Value = "'0.1"
$ws.Range("J4").Value = "calibMagnet.dat"
$ws.Range("K4").Value = ""
$ws.Range("M4").Value = "LockIn"
$ws.Range("N4").Value = "'0.0"

# --- Row 5 ---
$ws.Range("A5").Value = "100nm Py 4x3mm Insitu"
$ws.Range("B5").Value = "FieldSweep_27-10-22_15-31-14.dat"
$ws.Range("C5").Value = "100nm Py 4x3mm Insitu"
$ws.Range("D5").Value = "'15.0"
$ws.Range("E5").Value = "'13.0"
$ws.Range("F5").Value = "150/250/0.4"
$ws.Range("G5").Value = "'468.0"
$ws.Range("H5").Value = "'2.0"
$ws.Range("I5").Value = "'0.1"
$ws.Range("J5").Value = "calibMagnet.dat"
$ws.Range("K5").Value = "Changed Mod Freq"
$ws.Range("L5").Value = "In-Situ short"
$ws.Range("M5").Value = "LockIn"
$ws.Range("N5").Value = "'0.0"

# --- Row 6 ---
$ws.Range("A6").Value = "100nm Py 4x3mm Insitu"
$ws.Range("B6").Value = "FieldSweep_27-10-22_14-52-21.dat"
$ws.Range("C6").Value = "100nm Py 4x3mm Insitu"
$ws.Range("D6").Value = "'15.0"
$ws.Range("E6").Value = "'13.0"
$ws.Range("F6").Value = "150/250/0.4"
$ws.Range("G6").Value = "'3000.0"
$ws.Range("H6").Value = "'5.0"
$ws.Range("I6").Value = "'0.1"
$ws.Range("J6").Value = "calibMagnet.dat"
$ws.Range("K6").Value = ""
$ws.Range("L6").Value = "In-Situ short"
$ws.Range("M6").Value = "LockIn"
$ws.Range("N6").Value = "'0.0"

# --- Row 7 ---
$ws.Range("A7").Value = "100nm Py 4x3mm Insitu"
$ws.Range("B7").Value = "FieldSweep_27-10-22_14-30-06.dat"
$ws.Range("C7").Value = "100nm Py 4x3mm Insitu"
$ws.Range("D7").Value = "'15.0"
$ws.Range("E7").Value = "'13.0"
$ws.Range("F7").Value = "150/250/0.4"
$ws.Range("G7").Value = "'3000.0"
$ws.Range("H7").Value = "'5.0"
$ws.Range("I7").Value = "'0.1"
$ws.Range("J7").Value = "calibMagnet.dat"
$ws.Range("K7").Value = "Forgot to plug ion BNC cables"
$ws.Range("L7").Value = "In-Situ Short"
$ws.Range("M7").Value = "LockIn"
$ws.Range("N7").Value = "'0.0"

# --- Row 8 ---
$ws.Range("A8").Value = "TEsting Excel2"
$ws.Range("B8").Value = "Test"
$ws.Range("C8").Value = "None"
$ws.Range("D8").Value = "'15.0"
$ws.Range("E8").Value = "'13.0"
$ws.Range("F8").Value = "0/150/0.4"
$ws.Range("G8").Value = "'3000.0"
$ws.Range("H8").Value = "'5.0"
$ws.Range("I8").Value = "'0.1"
$ws.Range("J8").Value = "calibMagnet.dat"
$ws.Range("K8").Value = ""
$ws.Range("L8").Value = "/"
$ws.Range("M8").Value = "LockIn"
$ws.Range("N8").Value = "'0.0"

# --- Row 9 ---
$ws.Range("A9").Value = "Excel 1"
$ws.Range("B9").Value = "None"
$ws.Range("C9").Value = "None"
$ws.Range("D9").Value = "'15.0"
$ws.Range("E9").Value = "'13.0"
$ws.Range("F9").Value = "<PyQt5.QtWidgets.QLineEdit object at 0x000001237F4248B0>/<PyQt5.QtWidgets.QLineEdit object at 0x000001237F424940>/0.4"
$ws.Range("G9").Value = "'3000.0"
$ws.Range("H9").Value = "'5.0"
$ws.Range("I9").Value = "'0.1"
$ws.Range("J9").Value = "calibMagnet.dat"
$ws.Range("K9").Value = ""
$ws.Range("L9").Value = "/"
$ws.Range("M9").Value = "LockIn"
$ws.Range("N9").Value = "'0.0"

# --- Row 10 ---
$ws.Range("A10").Value = "/"
$ws.Range("B10").Value = "/"
$ws.Range("C10").Value = "/"
$ws.Range("D10").Value = "/"
$ws.Range("E10").Value = "/"
$ws.Range("F10").Value = "/"
$ws.Range("G10").Value = "/"
$ws.Range("H10").Value = "/"
$ws.Range("I10").Value = "/"
$ws.Range("J10").Value = "/"
$ws.Range("K10").Value = "/"
$ws.Range("L10").Value = "/"
$ws.Range("M10").Value = "/"
$ws.Range("N10").Value = "/"

# --- Row 11 ---
$ws.Range("A11").Value = "Test"
$ws.Range("B11").Value = "/"
$ws.Range("C11").Value = "/"
$ws.Range("D11").Value = "/"
$ws.Range("E11").Value = "Test"
$ws.Range("F11").Value = "/"
$ws.Range("G11").Value = "/"
$ws.Range("H11").Value = "/"
$ws.Range("I11").Value = "/"
$ws.Range("J11").Value = "/"
$ws.Range("K11").Value = "/"
$ws.Range("L11").Value = "/"
$ws.Range("M11").Value = "/"
$ws.Range("N11").Value = "/"

# --- Row 12 ---
$ws.Range("A12").Value = "/"
$ws.Range("B12").Value = "/"
$ws.Range("C12").Value = "/"
$ws.Range("D12").Value = "/"
$ws.Range("E12").Value = "/"
$ws.Range("F12").Value = "/"
$ws.Range("G12").Value = "/"
$ws.Range("H12").Value = "/"
$ws.Range("I12").Value = "/"
$ws.Range("J12").Value = "/"
$ws.Range("K12").Value = "/"
$ws.Range("L12").Value = "/"
$ws.Range("M12").Value = "/"
$ws.Range("N12").Value = "/"

